$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6: Days of Chunder / Antidote
$ws.Range("H6").Value = 166.5
$ws.Range("I6").Value = 166.5
$ws.Range("K6").Value = 499.5
$ws.Range("M6").Value = -387.5

# ALC row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 4708.9653
$ws.Range("I98").Value = 2849.2173
$ws.Range("J98").Value = 11838
$ws.Range("K98").Value = 2849.2173
$ws.Range("L98").Value = 11838
$ws.Range("M98").Value = -1351.2173
$ws.Range("N98").Value = -14834

# ALC row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1994.4688
$ws.Range("J112").Value = 2029.7742
$ws.Range("L112").Value = 6089.3226
$ws.Range("N112").Value = -8305.3226

# ALC row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 2671.923
$ws.Range("I113").Value = 2162
$ws.Range("J113").Value = 2990.625
$ws.Range("K113").Value = 2162
$ws.Range("L113").Value = 2990.625
$ws.Range("M113").Value = 1092
$ws.Range("N113").Value = -9498.625

# ALC row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 4708.9653
$ws.Range("I122").Value = 2849.2173
$ws.Range("J122").Value = 11838
$ws.Range("K122").Value = 8547.651899999999
$ws.Range("L122").Value = 35514
$ws.Range("M122").Value = -6097.651899999999
$ws.Range("N122").Value = -40414

# ALC row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 1199.5454
$ws.Range("I135").Value = 854
$ws.Range("J135").Value = 2374.4
$ws.Range("K135").Value = 7686
$ws.Range("L135").Value = 21369.6
$ws.Range("M135").Value = -5151
$ws.Range("N135").Value = -26439.6

# ALC row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 1540.75
$ws.Range("I137").Value = 1563.4286
$ws.Range("J137").Value = 1515.6842
$ws.Range("K137").Value = 4690.2858
$ws.Range("L137").Value = 4547.0526
$ws.Range("M137").Value = -2140.2858
$ws.Range("N137").Value = -9647.052599999999

# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1966.09
$ws.Range("I138").Value = 1203.2106
$ws.Range("J138").Value = 2145.037
$ws.Range("K138").Value = 3609.6318
$ws.Range("L138").Value = 6435.110999999999
$ws.Range("M138").Value = 1530.3682
$ws.Range("N138").Value = -16715.111

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 1525.7241
$ws.Range("I2").Value = 1093.1111
$ws.Range("J2").Value = 2233.6365
$ws.Range("K2").Value = 1093.1111
$ws.Range("L2").Value = 2233.6365
$ws.Range("M2").Value = -980.1111000000001
$ws.Range("N2").Value = -2459.6365

# ARM row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 2346.2856
$ws.Range("I45").Value = 2084.8
$ws.Range("K45").Value = 2084.8
$ws.Range("M45").Value = -1707.8

# ARM row 108: Time to Fry / Deepgold Rail Frypan
$ws.Range("H108").Value = 46000
$ws.Range("J108").Value = 46000
$ws.Range("L108").Value = 46000
$ws.Range("N108").Value = -53680

# ARM row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 1138.6364
$ws.Range("I110").Value = 1115.8
$ws.Range("J110").Value = 1187.5714
$ws.Range("K110").Value = 1115.8
$ws.Range("L110").Value = 1187.5714
$ws.Range("M110").Value = 929.2
$ws.Range("N110").Value = -5277.5714

# ARM row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 1525.7241
$ws.Range("I116").Value = 1093.1111
$ws.Range("J116").Value = 2233.6365
$ws.Range("K116").Value = 1093.1111
$ws.Range("L116").Value = 2233.6365
$ws.Range("M116").Value = 1200.8889
$ws.Range("N116").Value = -6821.636500000001

# ARM row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2144.2222
$ws.Range("I122").Value = 1916
$ws.Range("K122").Value = 5748
$ws.Range("M122").Value = -3298

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 1525.7241
$ws.Range("I3").Value = 1093.1111
$ws.Range("J3").Value = 2233.6365
$ws.Range("K3").Value = 1093.1111
$ws.Range("L3").Value = 2233.6365
$ws.Range("M3").Value = -979.1111000000001
$ws.Range("N3").Value = -2461.6365

# BSM row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2093.923
$ws.Range("I134").Value = 1856.2572
$ws.Range("J134").Value = 4173.5
$ws.Range("K134").Value = 5568.7716
$ws.Range("L134").Value = 12520.5
$ws.Range("M134").Value = -3033.7716
$ws.Range("N134").Value = -17590.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 1666.3334
$ws.Range("I16").Value = 1666.3334
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1666.3334
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1379.3334
$ws.Range("N16").ClearContents()

# CRP row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 1666.3334
$ws.Range("I113").Value = 1666.3334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1666.3334
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 503.6666
$ws.Range("N113").ClearContents()

# CRP row 114: Ground to a Halt / White Ash Grinding Wheel
$ws.Range("H114").Value = 39799.5
$ws.Range("I114").Value = 35000
$ws.Range("J114").Value = 44599
$ws.Range("K114").Value = 35000
$ws.Range("L114").Value = 44599
$ws.Range("M114").Value = -30661
$ws.Range("N114").Value = -53277

$ws = $wb.Worksheets.Item("CUL")
# CUL row 64: The Aroma of Faith / Baked Onion Soup
$ws.Range("H64").Value = 6571.2856
$ws.Range("J64").Value = 6999.8335
$ws.Range("L64").Value = 20999.5005
$ws.Range("N64").Value = -21539.5005

# CUL row 67: Soup's On (L) / Baked Onion Soup
$ws.Range("H67").Value = 6571.2856
$ws.Range("J67").Value = 6999.8335
$ws.Range("L67").Value = 20999.5005
$ws.Range("N67").Value = -22871.5005

# CUL row 104: Fits to a Tea / Doman Tea
$ws.Range("H104").Value = 7500
$ws.Range("J104").Value = 7500
$ws.Range("L104").Value = 22500
$ws.Range("N104").Value = -27742

# CUL row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 10871902
$ws.Range("J131").Value = 11629217
$ws.Range("L131").Value = 34887651
$ws.Range("N131").Value = -34897731

# CUL row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 33335866
$ws.Range("I137").Value = 1073.3334
$ws.Range("J137").Value = 47622204
$ws.Range("K137").Value = 3220.0002
$ws.Range("L137").Value = 142866612
$ws.Range("M137").Value = 1879.9998
$ws.Range("N137").Value = -142876812

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 3165.5417
$ws.Range("I122").Value = 3227.762
$ws.Range("J122").Value = 2730
$ws.Range("K122").Value = 9683.286
$ws.Range("L122").Value = 8190
$ws.Range("M122").Value = -7233.286
$ws.Range("N122").Value = -13090

# GSM row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 2706
$ws.Range("I126").Value = 2509.6
$ws.Range("J126").Value = 2902.4
$ws.Range("K126").Value = 7528.799999999999
$ws.Range("L126").Value = 8707.200000000001
$ws.Range("M126").Value = -5058.799999999999
$ws.Range("N126").Value = -13647.2

$ws = $wb.Worksheets.Item("LTW")
# LTW row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1994.4445
$ws.Range("I61").Value = 1900.0667
$ws.Range("J61").Value = 2466.3333
$ws.Range("K61").Value = 1900.0667
$ws.Range("L61").Value = 2466.3333
$ws.Range("M61").Value = -1698.0667
$ws.Range("N61").Value = -2870.3333

# LTW row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1994.4445
$ws.Range("I113").Value = 1900.0667
$ws.Range("J113").Value = 2466.3333
$ws.Range("K113").Value = 1900.0667
$ws.Range("L113").Value = 2466.3333
$ws.Range("M113").Value = 269.9332999999999
$ws.Range("N113").Value = -6806.3333

# LTW row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 22733090
$ws.Range("I122").Value = 27783776
$ws.Range("K122").Value = 83351328
$ws.Range("M122").Value = -83348878

$ws = $wb.Worksheets.Item("WVR")
# WVR row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 734.62067
$ws.Range("J113").Value = 1664
$ws.Range("L113").Value = 4992
$ws.Range("N113").Value = -9332

# WVR row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 7793.8184
$ws.Range("I126").Value = 8472.799999999999
$ws.Range("J126").Value = 1004
$ws.Range("K126").Value = 25418.4
$ws.Range("L126").Value = 3012
$ws.Range("M126").Value = -22948.4
$ws.Range("N126").Value = -7952
